# Grade tables now also in stripped form for easier input.
# Fill column A (rows 10-33) with "X", matching the existing A9 marker,
# and reset the active selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 10; $row -le 33; $row++) {
    $ws.Cells.Item($row, 1).Value = "X"
}

$ws.Range("A1").Select()
